$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Option 1 - LR1 - DN1 (80-20)")
$ws.Range("D7").Value = 0.0179918809253119
$ws.Range("E7").Value = 0.1186304499816319
$ws.Range("F7").Value = 0.1696272826101057
$ws.Range("G7").Value = 0.1341338172323143
$ws.Range("H7").Value = 18.00197642917625
$ws.Range("L7").Value = 0.03250270772802855
$ws.Range("M7").Value = 0.1480394330773291
$ws.Range("N7").Value = 0.204227855967163
$ws.Range("O7").Value = 0.1802850735031288
$ws.Range("P7").Value = 23.89406442753663

$ws = $wb.Worksheets.Item("Option 1 - LR1 - DN1 (60-40)")
$ws.Range("D7").Value = 0.0281175482060782
$ws.Range("E7").Value = 0.1279678585492715
$ws.Range("F7").Value = 0.1902326352532762
$ws.Range("G7").Value = 0.1676828798836608
$ws.Range("H7").Value = 20.21080586637246
$ws.Range("L7").Value = 0.02466042671762571
$ws.Range("M7").Value = 0.1113212340714586
$ws.Range("N7").Value = 0.171858315769591
$ws.Range("O7").Value = 0.1570363866039515
$ws.Range("P7").Value = 17.67252413352941

$ws = $wb.Worksheets.Item("Option 1 - LR1 - DN1 (70-30)")
$ws.Range("D7").Value = 0.01287121123128601
$ws.Range("E7").Value = 0.08563202834452491
$ws.Range("F7").Value = 0.1181390320346846
$ws.Range("G7").Value = 0.1134513606409637
$ws.Range("H7").Value = 12.53255944509988
$ws.Range("L7").Value = 0.01998722168110157
$ws.Range("M7").Value = 0.1209672007108275
$ws.Range("N7").Value = 0.1790174201338875
$ws.Range("O7").Value = 0.1413761708390122
$ws.Range("P7").Value = 19.98701887951244

$ws = $wb.Worksheets.Item("Option 1 - LR1 - DN2 (80-20)")
$ws.Range("D7").Value = 0.08879198112204009
$ws.Range("E7").Value = 0.2651541837938986
$ws.Range("F7").Value = 0.09584498273079374
$ws.Range("G7").Value = 0.297979833415015
$ws.Range("H7").Value = 9.926561799027152
$ws.Range("L7").Value = 0.1329966971440318
$ws.Range("M7").Value = 0.2993721930238203
$ws.Range("N7").Value = 0.1081281580231559
$ws.Range("O7").Value = 0.3646871222624015
$ws.Range("P7").Value = 11.74463826906958

$ws = $wb.Worksheets.Item("Option 1 - LR1 - DN2 (60-40)")
$ws.Range("D7").Value = 0.1049348382747722
$ws.Range("E7").Value = 0.2151390074506735
$ws.Range("F7").Value = 0.09334073671045441
$ws.Range("G7").Value = 0.323936472591112
$ws.Range("H7").Value = 8.909510778174404
$ws.Range("L7").Value = 0.1319980439940538
$ws.Range("M7").Value = 0.2680672671904306
$ws.Range("N7").Value = 0.1062093585897546
$ws.Range("O7").Value = 0.3633153506171379
$ws.Range("P7").Value = 10.71208133155415

$ws = $wb.Worksheets.Item("Option 1 - LR1 - DN2 (70-30)")
$ws.Range("D7").Value = 0.07472352836974828
$ws.Range("E7").Value = 0.1919161812161036
$ws.Range("F7").Value = 0.06744969770511261
$ws.Range("G7").Value = 0.2733560468871107
$ws.Range("H7").Value = 7.178312920982592
$ws.Range("L7").Value = 0.1029953252738324
$ws.Range("M7").Value = 0.2359410632715214
$ws.Range("N7").Value = 0.08245310106045689
$ws.Range("O7").Value = 0.3209288476809656
$ws.Range("P7").Value = 8.874698031978594

$ws = $wb.Worksheets.Item("Option 1 - LR2 - DN1 (80-20)")
$ws.Range("D7").Value = 1.056994605137028
$ws.Range("E7").Value = 0.858724435680976
$ws.Range("F7").Value = 883717275121979.8
$ws.Range("G7").Value = 1.028102429302172
$ws.Range("H7").Value = 177.0027571523113
$ws.Range("L7").Value = 1.033145897545401
$ws.Range("M7").Value = 0.8566940721454438
$ws.Range("N7").Value = 831400841312325.9
$ws.Range("O7").Value = 1.016437847359789
$ws.Range("P7").Value = 180.5183506677184

$ws = $wb.Worksheets.Item("Option 1 - LR2 - DN1 (60-40)")
$ws.Range("D7").Value = 0.4051472972625934
$ws.Range("E7").Value = 0.4894173230022899
$ws.Range("F7").Value = 1072036993429984
$ws.Range("G7").Value = 0.6365118202065013
$ws.Range("H7").Value = 83.26798026807425
$ws.Range("L7").Value = 0.4176178695573639
$ws.Range("M7").Value = 0.5207142485266109
$ws.Range("N7").Value = 1196903332573756
$ws.Range("O7").Value = 0.6462336029311412
$ws.Range("P7").Value = 81.66427961802161

$ws = $wb.Worksheets.Item("Option 1 - LR2 - DN1 (70-30)")
$ws.Range("D7").Value = 1.016765273733438
$ws.Range("E7").Value = 0.8550418001520068
$ws.Range("F7").Value = 870896314341999
$ws.Range("G7").Value = 1.008347794034101
$ws.Range("H7").Value = 159.8109631189442
$ws.Range("L7").Value = 1.203211406246794
$ws.Range("M7").Value = 0.9222647529419414
$ws.Range("N7").Value = 826020043883285
$ws.Range("O7").Value = 1.096909935339631
$ws.Range("P7").Value = 167.2729577605253

$ws = $wb.Worksheets.Item("Option 1 - LR2 - DN2 (60-40)")
$ws.Range("D7").Value = 1.113949853911995
$ws.Range("E7").Value = 0.8811127075855789
$ws.Range("F7").Value = 1239582830352277
$ws.Range("G7").Value = 1.055438228373407
$ws.Range("H7").Value = 130.7315789378721
$ws.Range("L7").Value = 1.970900756794443
$ws.Range("M7").Value = 1.1365416461419
$ws.Range("N7").Value = 1878223606308148
$ws.Range("O7").Value = 1.403887729412307
$ws.Range("P7").Value = 127.2944179994935

$ws = $wb.Worksheets.Item("Option 1 - LR2 - DN2 (70-30)")
$ws.Range("D7").Value = 2.626297509864584
$ws.Range("E7").Value = 1.307411706232485
$ws.Range("F7").Value = 1924762208522051
$ws.Range("G7").Value = 1.620585545370742
$ws.Range("H7").Value = 147.5295062912983
$ws.Range("L7").Value = 8.961966498157992
$ws.Range("M7").Value = 2.399776615514196
$ws.Range("N7").Value = 3207652147597913
$ws.Range("O7").Value = 2.99365437186025
$ws.Range("P7").Value = 164.287008092179

$ws = $wb.Worksheets.Item("Option 1 - LR2 - DN2 (80-20)")
$ws.Range("D7").Value = 2.187003232381139
$ws.Range("E7").Value = 1.260030062154847
$ws.Range("F7").Value = 1372020673264916
$ws.Range("G7").Value = 1.47885199813272
$ws.Range("H7").Value = 192.6092387796874
$ws.Range("L7").Value = 3.749116858523515
$ws.Range("M7").Value = 1.618175000014997
$ws.Range("N7").Value = 2895850310827128
$ws.Range("O7").Value = 1.936263633528119
$ws.Range("P7").Value = 182.2623187051494

$ws = $wb.Worksheets.Item("Option 1 - NLR1 - DN1 (80-20)")
$ws.Range("D7").Value = 0.04905374744615793
$ws.Range("E7").Value = 0.1348942525282142
$ws.Range("F7").Value = 1.180031753678706
$ws.Range("G7").Value = 0.2214808060445824
$ws.Range("H7").Value = 27.47842772577463
$ws.Range("L7").Value = 0.04854120948959907
$ws.Range("M7").Value = 0.1344970376261711
$ws.Range("N7").Value = 1.17341376202927
$ws.Range("O7").Value = 0.220320696916107
$ws.Range("P7").Value = 27.49387665851739

$ws = $wb.Worksheets.Item("Option 1 - NLR1 - DN1 (70-30)")
$ws.Range("D7").Value = 0.03076230876773942
$ws.Range("E7").Value = 0.10065807139648
$ws.Range("F7").Value = 1.061274415764125
$ws.Range("G7").Value = 0.1753918720116169
$ws.Range("H7").Value = 22.56021813802044
$ws.Range("L7").Value = 0.03082799117911534
$ws.Range("M7").Value = 0.1006783586476175
$ws.Range("N7").Value = 1.063519356121315
$ws.Range("O7").Value = 0.1755790169100948
$ws.Range("P7").Value = 22.56724825816929

$ws = $wb.Worksheets.Item("Option 1 - NLR1 - DN1 (60-40)")
$ws.Range("D7").Value = 0.03145118012700405
$ws.Range("E7").Value = 0.1176109511430981
$ws.Range("F7").Value = 1.064104257804336
$ws.Range("G7").Value = 0.1773448057514063
$ws.Range("H7").Value = 28.0377848840399
$ws.Range("L7").Value = 0.02877552441652384
$ws.Range("M7").Value = 0.1025793405321778
$ws.Range("N7").Value = 1.05677053477698
$ws.Range("O7").Value = 0.169633500277875
$ws.Range("P7").Value = 23.62401367184036

$ws = $wb.Worksheets.Item("Option 1 - NLR1 - DN2 (80-20)")
$ws.Range("D7").Value = 0.5211246120532131
$ws.Range("E7").Value = 0.3336833755184812
$ws.Range("F7").Value = 0.5727384682036503
$ws.Range("G7").Value = 0.7218896120967617
$ws.Range("H7").Value = 16.81272907165241
$ws.Range("L7").Value = 0.5209976961238756
$ws.Range("M7").Value = 0.3336778446090488
$ws.Range("N7").Value = 0.5726927452202003
$ws.Range("O7").Value = 0.7218017013861048
$ws.Range("P7").Value = 16.8131097436111

$ws = $wb.Worksheets.Item("Option 1 - NLR1 - DN2 (70-30)")
$ws.Range("D7").Value = 0.3963140912831225
$ws.Range("E7").Value = 0.269534297037626
$ws.Range("F7").Value = 0.5259475840008978
$ws.Range("G7").Value = 0.629534821342809
$ws.Range("H7").Value = 14.41449720095513
$ws.Range("L7").Value = 0.4397589409703375
$ws.Range("M7").Value = 0.2802417544319984
$ws.Range("N7").Value = 0.551743323697844
$ws.Range("O7").Value = 0.6631432280965684
$ws.Range("P7").Value = 14.86114731823628

$ws = $wb.Worksheets.Item("Option 1 - NLR1 - DN2 (60-40)")
$ws.Range("D7").Value = 0.3341419342718551
$ws.Range("E7").Value = 0.2808505500743309
$ws.Range("F7").Value = 0.2446381002270317
$ws.Range("G7").Value = 0.5780501139796229
$ws.Range("H7").Value = 14.38313770986017
$ws.Range("L7").Value = 0.340126881879963
$ws.Range("M7").Value = 0.290135949221715
$ws.Range("N7").Value = 0.2495320181462332
$ws.Range("O7").Value = 0.5832039796503132
$ws.Range("P7").Value = 14.79341989135393

$ws = $wb.Worksheets.Item("Option 1 - NLR2 - DN1 (80-20)")
$ws.Range("D7").Value = 0.04024338070258342
$ws.Range("E7").Value = 0.1568905693894855
$ws.Range("F7").Value = 0.2316665929312102
$ws.Range("G7").Value = 0.2006075290276599
$ws.Range("H7").Value = 21.58477685365705

$ws = $wb.Worksheets.Item("Option 1 - NLR2 - DN1 (70-30)")
$ws.Range("D7").Value = 0.1713580497920794
$ws.Range("E7").Value = 0.3414249082301279
$ws.Range("F7").Value = 0.5415784986528271
$ws.Range("G7").Value = 0.4139541638781755
$ws.Range("H7").Value = 68.92945276192613

$ws = $wb.Worksheets.Item("Option 1 - NLR2 - DN1 (60-40)")
$ws.Range("D7").Value = 2.015034735724734
$ws.Range("E7").Value = 1.017926757261285
$ws.Range("F7").Value = 2.42821152037283
$ws.Range("G7").Value = 1.419519191742307
$ws.Range("H7").Value = 114.935707975681

$ws = $wb.Worksheets.Item("Option 1 - NLR2 - DN2 (80-20)")
$ws.Range("D7").Value = 0.9417912160220516
$ws.Range("E7").Value = 0.7595408877367373
$ws.Range("F7").Value = 0.2704198110401613
$ws.Range("G7").Value = 0.97045928097064
$ws.Range("H7").Value = 25.44736189641421

$ws = $wb.Worksheets.Item("Option 1 - NLR2 - DN2 (70-30)")
$ws.Range("D7").Value = 4.748770044021249
$ws.Range("E7").Value = 1.75983036844142
$ws.Range("F7").Value = 0.6583525225842659
$ws.Range("G7").Value = 2.179167282248256
$ws.Range("H7").Value = 84.51012561520299

$ws = $wb.Worksheets.Item("Option 1 - NLR2 - DN2 (60-40)")
$ws.Range("D7").Value = 1.145750813661144
$ws.Range("E7").Value = 0.8071426284398897
$ws.Range("F7").Value = 0.3005445726125022
$ws.Range("G7").Value = 1.070397502641492
$ws.Range("H7").Value = 33.14846474252346
